$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# The hyperlink currently anchored on B27 (-> SNCE.xlsx) needs to end up on
# B28 once the new row is inserted below it. This engine does not shift
# hyperlink anchors automatically the way Excel does on a real row insert,
# and per-link deletion is a no-op here, so the whole-sheet collection has
# to be cleared (the only reliable deletion API available) and every link
# rebuilt afterwards at its correct final address.
$ws.Range("B3").Hyperlinks.Delete()

# Insert a new row at row 18, pushing old rows 18+ down by one.
$ws.Rows.Item(18).Insert()

# Row 17 (Amerisource / ABC) now also gets the "x" flag in column A.
$ws.Range("A17").Value = "x"

# New row 18: JD Health / 6618 HK, flagged with "x" in column A.
$ws.Range("A18").Value = "x"
$ws.Range("B18").Value = "JD Health"
$ws.Range("C18").Value = "6618 HK"

# Rebuild the three hyperlinks at their final addresses, and restore the
# original cell styles afterwards so re-adding the links doesn't silently
# recolor/re-underline those cells (Hyperlinks.Add reapplies its own
# "Hyperlink" character style on top of whatever style the cell already has).
$b3Style = $ws.Range("B3").Style
$b4Style = $ws.Range("B4").Style
$b28Style = $ws.Range("B28").Style

$ws.Hyperlinks.Add($ws.Range("B28"), "SNCE.xlsx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "UNH.xlsx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "CVS.xlsx") | Out-Null

$ws.Range("B3").Style = $b3Style
$ws.Range("B4").Style = $b4Style
$ws.Range("B28").Style = $b28Style

# Update selection to match the saved view state.
$ws.Range("A19").Select()
